# Fixing/cleaning the "Length group in feet" labels in column A of the
# nvessels table (Table5) and restoring a sane ascending ordering of the
# groups. Also refreshes the active view (scrolled / selection) to match
# where the author ended up after the cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list (not a hashtable) so the cells get rewritten in the exact
# sequence the author typed them in, matching the resulting shared-string
# table layout.
$labels = @(
    @{ Row = 21; Text = '101-105' },
    @{ Row = 23; Text = '111-115' },
    @{ Row = 24; Text = '116-120' },
    @{ Row = 25; Text = '121-125' },
    @{ Row = 27; Text = '131-135' },
    @{ Row = 29; Text = '141-145' },
    @{ Row = 30; Text = '146-150' },
    @{ Row = 31; Text = '151-155' },
    @{ Row = 35; Text = '171-175' },
    @{ Row = 36; Text = '176-180' },
    @{ Row = 38; Text = 'Total' },
    @{ Row = 4;  Text = '16- 20' },
    @{ Row = 6;  Text = '26- 30' },
    @{ Row = 8;  Text = '36- 40 ' },
    @{ Row = 10; Text = '46- 50 ' },
    @{ Row = 12; Text = '56- 60' },
    @{ Row = 13; Text = '61- 65' },
    @{ Row = 14; Text = '66- 70' },
    @{ Row = 16; Text = '76- 80' },
    @{ Row = 18; Text = '86- 90' },
    @{ Row = 19; Text = '91- 95' },
    @{ Row = 33; Text = '161-165' },
    @{ Row = 20; Text = '96- 100' },
    @{ Row = 22; Text = '106-110' },
    @{ Row = 26; Text = '126-130' },
    @{ Row = 28; Text = '136-140' },
    @{ Row = 32; Text = '156-100' },
    @{ Row = 34; Text = '166-170' }
)

foreach ($item in $labels) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Text
}

# Scroll the view down and move the selection, mirroring where the author
# left the sheet after finishing the cleanup.
$ws.Range("A40").Select()
$excel.ActiveWindow.ScrollRow = 18
